$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# tc0011 now also covers applyKcatConstraints (row 4) and selectKCats (row 13),
# and partially covers writeDLKcatInput (row 14) and mergeDlkcatAndFuzzyKcats (row 10).
$ws.Range("B4").Value = "tc0011"
$ws.Range("B10").Value = "tc0011"
$ws.Range("B14").Value = "partly by tc0011"
$ws.Range("B13").Value = "tc0011 - we did not test all possible parameterizations here"

# Move the active selection from B9 to B17
$ws.Range("B17").Select() | Out-Null
